$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 25004304
$ws.Range("J32").Value = 4281.125
$ws.Range("L32").Value = 4281.125
$ws.Range("N32").Value = -4933.125
$ws.Range("H33").Value = 153.45
$ws.Range("I33").Value = 178.86667
$ws.Range("J33").Value = 77.2
$ws.Range("K33").Value = 178.86667
$ws.Range("L33").Value = 77.2
$ws.Range("M33").Value = 50.13333
$ws.Range("N33").Value = -535.2
$ws.Range("H40").Value = 71468290
$ws.Range("I40").Value = 46332.332
$ws.Range("K40").Value = 46332.332
$ws.Range("M40").Value = -46157.332
$ws.Range("H62").Value = 3133.6667
$ws.Range("I62").Value = 2901
$ws.Range("K62").Value = 2901
$ws.Range("M62").Value = -2277
$ws.Range("H65").Value = 3133.6667
$ws.Range("I65").Value = 2901
$ws.Range("K65").Value = 14505
$ws.Range("M65").Value = -11385
$ws.Range("H106").Value = 17546170
$ws.Range("I106").Value = 22223780
$ws.Range("J106").Value = 5126.25
$ws.Range("K106").Value = 22223780
$ws.Range("L106").Value = 5126.25
$ws.Range("M106").Value = -22223149
$ws.Range("N106").Value = -6388.25
$ws.Range("H125").Value = 2298.6
$ws.Range("I125").Value = 864.6667
$ws.Range("J125").Value = 4449.5
$ws.Range("K125").Value = 7782.0003
$ws.Range("L125").Value = 40045.5
$ws.Range("M125").Value = -5322.0003
$ws.Range("N125").Value = -44965.5
$ws.Range("H137").Value = 3550.6316
$ws.Range("I137").Value = 1523.625
$ws.Range("K137").Value = 4570.875
$ws.Range("M137").Value = -2020.875
$ws.Range("H141").Value = 6833.5625
$ws.Range("I141").Value = 5667.2856
$ws.Range("K141").Value = 17001.8568
$ws.Range("M141").Value = -11821.8568

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3414.44
$ws.Range("I32").Value = 1598.3096
$ws.Range("K32").Value = 1598.3096
$ws.Range("M32").Value = -1311.3096
$ws.Range("H45").Value = 1807.7
$ws.Range("I45").Value = 1826
$ws.Range("K45").Value = 1826
$ws.Range("M45").Value = -1449
$ws.Range("H88").Value = 101280
$ws.Range("J88").Value = 126175
$ws.Range("L88").Value = 126175
$ws.Range("N88").Value = -126987
$ws.Range("H91").Value = 101280
$ws.Range("J91").Value = 126175
$ws.Range("L91").Value = 126175
$ws.Range("N91").Value = -128983
$ws.Range("H132").Value = 30787.87
$ws.Range("I132").Value = 45001.77
$ws.Range("K132").Value = 135005.31
$ws.Range("M132").Value = -132475.31

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2478.7827
$ws.Range("I31").Value = 963.375
$ws.Range("J31").Value = 5942.5713
$ws.Range("K31").Value = 963.375
$ws.Range("L31").Value = 5942.5713
$ws.Range("M31").Value = -668.375
$ws.Range("N31").Value = -6532.5713
$ws.Range("H34").Value = 2478.7827
$ws.Range("I34").Value = 963.375
$ws.Range("J34").Value = 5942.5713
$ws.Range("K34").Value = 963.375
$ws.Range("L34").Value = 5942.5713
$ws.Range("M34").Value = -761.375
$ws.Range("N34").Value = -6346.5713
$ws.Range("H60").Value = 234333
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H105").Value = 1624261.1
$ws.Range("I105").Value = 1894637.9
$ws.Range("K105").Value = 1894637.9
$ws.Range("M105").Value = -1892890.9

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 869.25
$ws.Range("I113").Value = 759
$ws.Range("K113").Value = 2277
$ws.Range("M113").Value = -107
$ws.Range("H122").Value = 494.33334
$ws.Range("I122").Value = 496.66666
$ws.Range("K122").Value = 4469.99994
$ws.Range("M122").Value = -2019.99994
$ws.Range("H129").Value = 2181
$ws.Range("I129").Value = 1097.4286
$ws.Range("J129").Value = 3129.125
$ws.Range("K129").Value = 3292.2858
$ws.Range("L129").Value = 9387.375
$ws.Range("M129").Value = 1707.7142
$ws.Range("N129").Value = -19387.375
$ws.Range("H131").Value = 2565.7334
$ws.Range("J131").Value = 2761.585
$ws.Range("L131").Value = 8284.755000000001
$ws.Range("N131").Value = -18364.755
$ws.Range("H138").Value = 66195.94
$ws.Range("I138").Value = 128475
$ws.Range("J138").Value = 3916.875
$ws.Range("K138").Value = 385425
$ws.Range("L138").Value = 11750.625
$ws.Range("M138").Value = -380285
$ws.Range("N138").Value = -22030.625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 13458.333
$ws.Range("I102").Value = 14611.111
$ws.Range("K102").Value = 14611.111
$ws.Range("M102").Value = -12989.111
$ws.Range("H122").Value = 689284.9
$ws.Range("I122").Value = 1101699.1
$ws.Range("K122").Value = 3305097.3
$ws.Range("M122").Value = -3302647.3
$ws.Range("H132").Value = 4158.727
$ws.Range("I132").Value = 3741.5
$ws.Range("K132").Value = 11224.5
$ws.Range("M132").Value = -8694.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5527.6
$ws.Range("I7").Value = 3966.56
$ws.Range("J7").Value = 8129.3335
$ws.Range("K7").Value = 3966.56
$ws.Range("L7").Value = 8129.3335
$ws.Range("M7").Value = -3854.56
$ws.Range("N7").Value = -8353.333500000001
$ws.Range("H40").Value = 7130.7144
$ws.Range("I40").Value = 6668
$ws.Range("K40").Value = 6668
$ws.Range("M40").Value = -6532
$ws.Range("H122").Value = 166675920
$ws.Range("I122").Value = 500013500
$ws.Range("J122").Value = 7136
$ws.Range("K122").Value = 1500040500
$ws.Range("L122").Value = 21408
$ws.Range("M122").Value = -1500038050
$ws.Range("N122").Value = -26308
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 5527.6
$ws.Range("I126").Value = 3966.56
$ws.Range("J126").Value = 8129.3335
$ws.Range("K126").Value = 11899.68
$ws.Range("L126").Value = 24388.0005
$ws.Range("M126").Value = -9429.68
$ws.Range("N126").Value = -29328.0005
$ws.Range("H132").Value = 4184.1816
$ws.Range("I132").Value = 3049.5
$ws.Range("K132").Value = 9148.5
$ws.Range("M132").Value = -6618.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value = 34000
$ws.Range("J59").Value = 34000
$ws.Range("L59").Value = 34000
$ws.Range("N59").Value = -35476
$ws.Range("H126").Value = 2124.25
$ws.Range("I126").Value = 1998.5
$ws.Range("J126").Value = 2250
$ws.Range("K126").Value = 5995.5
$ws.Range("L126").Value = 6750
$ws.Range("M126").Value = -3525.5
$ws.Range("N126").Value = -11690
$ws.Range("H132").Value = 16133314
$ws.Range("I132").Value = 1884.8695
$ws.Range("J132").Value = 62511172
$ws.Range("K132").Value = 5654.6085
$ws.Range("L132").Value = 187533516
$ws.Range("M132").Value = -3124.6085
$ws.Range("N132").Value = -187538576
